$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "REX_DEF" in column F, matching the style of the existing
# header cells (B1:E1), which use style index 1 (bold, centered, bordered).
$ws.Range("F1").Value = "REX_DEF"
$ws.Range("F1").Style = $ws.Range("E1").Style

# Also copy the full formatting (font, border, alignment) from E1 to F1 to be safe.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
